$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Remove the empty "Address" (F) cells that had no content (inline empty string)
$ws.Range("F4").ClearContents()
$ws.Range("F15").ClearContents()
$ws.Range("F20").ClearContents()
$ws.Range("F39").ClearContents()

# Correct district names to official names
$ws.Range("G16").Value = "Bagalkot"
$ws.Range("G19").Value = "Chikkamagaluru (Chikmagalur)"
$ws.Range("G31").Value = "Dharwad"
$ws.Range("G34").Value = "Kalaburagi (Gulbarga)"
$ws.Range("G41").Value = "Kalaburagi (Gulbarga)"
$ws.Range("G46").Value = "Vijayapura (Bijapur)"
$ws.Range("G56").Value = "Vijayapura (Bijapur)"
$ws.Range("G57").Value = "Davangere"
$ws.Range("G58").Value = "Davangere"
